# "Reversing to my last commit" - revert the Sprint1 Burndown chart sheet
# back to its prior state: rename the "Meetings" task back to
# "Meeting semanal" with reduced logged effort, and drop the
# "Fazer pdf dos 3 User Stories..." task (row 15) entirely, clearing the
# hours that had been logged against it / the other tasks on Day 6 (col K).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - "Jogar o jogo para conhecer melhor o projeto": initial estimate
# drops from 8 to 7, and the Day 6 (K) effort entry is removed.
$ws.Range("D6").Value = 7
$ws.Range("K6").ClearContents()

# Row 7 - "Dar ideias no servidor de discord e discuti-las": Day 6 (K)
# effort entry removed.
$ws.Range("K7").ClearContents()

# Row 11 - task renamed from "Meetings" back to "Meeting semanal", with a
# smaller initial estimate and no Day 6 effort logged.
$ws.Range("C11").Value = "Meeting semanal"
$ws.Range("D11").Value = 2
$ws.Range("K11").ClearContents()

# Row 13 - "Analisar o código dado": Day 4/5/6 (I,J,K) effort entries removed.
$ws.Range("I13").ClearContents()
$ws.Range("J13").ClearContents()
$ws.Range("K13").ClearContents()

# Row 14 - "Pesquisar como se joga o jogo": Day 6 (K) effort entry removed.
$ws.Range("K14").ClearContents()

# Row 15 - "Fazer pdf dos 3 User Stories mais votados pela equipa e
# submeter no moodle" task is removed entirely (description, estimate and
# logged effort all cleared).
$ws.Range("C15").ClearContents()
$ws.Range("D15").ClearContents()
$ws.Range("K15").ClearContents()

# Restore the previous selection/view state.
[void]$ws.Range("N12").Select()
